$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.308.77"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.45%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.595.98"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.37%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E6").Value = "  -0.16%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  -0.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.02"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.07%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0855"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.22%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.820.90"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.41%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.596.87"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.57%  "
$ws.Range("E14").Value = "  -0.73%  "
$ws.Range("E15").Value = "  -1.38%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.44"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.31%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.297.08"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.48%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "229.83"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +6.96%  "
$ws.Range("E19").Value = "  +3.83%  "
$ws.Range("E20").Value = "  -0.47%  "
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("E22").Value = "  -0.26%  "
$ws.Range("E23").Value = "  +2.38%  "
$ws.Range("E24").Value = "  -1.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.39"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.97"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.02%  "
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.40"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.15%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0493"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.06%  "
$ws.Range("E31").Value = "  -0.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.492.92"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.05%  "
$ws.Range("E33").Value = "  +0.87%  "
$ws.Range("E34").Value = "  -0.87%  "
$ws.Range("E35").Value = "  -0.40%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.47"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.26%  "
$ws.Range("E37").Value = "  -3.24%  "
$ws.Range("E38").Value = "  -0.44%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.816"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.00%  "
$ws.Range("E40").Value = "  -1.65%  "
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("E42").Value = "  +1.33%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.936"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.733.31"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.49%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.758"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.59%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "60.35"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.20%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "88.45"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.42%  "
$ws.Range("E48").Value = "  -0.49%  "
$ws.Range("E49").Value = "  -0.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0955"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.57%  "
$ws.Range("E51").Value = "  +0.05%  "
